$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.293.93"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.775.65"
$ws.Range("E3").Value = "  +3.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.43"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5142"
$ws.Range("E7").Value = "  +7.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3684"
$ws.Range("E8").Value = "  +6.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.74"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07388"
$ws.Range("E10").Value = "  +1.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.087"
$ws.Range("E11").Value = "  +3.88%  "

$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.47"
$ws.Range("E13").Value = "  +2.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.068"
$ws.Range("E14").Value = "  +3.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.767.46"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.953"
$ws.Range("E16").Value = "  +1.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.12"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001046"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06433"
$ws.Range("E19").Value = "  +1.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.75"
$ws.Range("E21").Value = "  +1.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.828"
$ws.Range("E22").Value = "  +3.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.328.18"
$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.23"
$ws.Range("E24").Value = "  +3.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.119"
$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.56"
$ws.Range("E26").Value = "  +1.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.19"
$ws.Range("E27").Value = "  +2.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.328"
$ws.Range("E28").Value = "  +10.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.974.19"
$ws.Range("E29").Value = "  +3.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.15"
$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.067"
$ws.Range("E31").Value = "  +4.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09745"
$ws.Range("E32").Value = "  +5.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.573"
$ws.Range("E33").Value = "  +4.42%  "

$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02243"
$ws.Range("E35").Value = "  +1.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05973"
$ws.Range("E36").Value = "  +0.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.24"
$ws.Range("E37").Value = "  +1.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6149"
$ws.Range("E38").Value = "  +3.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.838"
$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2020"
$ws.Range("E40").Value = "  +0.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.436"
$ws.Range("E41").Value = "  +1.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.090"
$ws.Range("E42").Value = "  +8.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.137"
$ws.Range("E43").Value = "  +3.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.14"
$ws.Range("E44").Value = "  +3.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5763"
$ws.Range("E45").Value = "  +2.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.630"
$ws.Range("E46").Value = "  +1.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.28"
$ws.Range("E47").Value = "  +2.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.885"
$ws.Range("E48").Value = "  +2.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.114"
$ws.Range("E49").Value = "  +2.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06706"
$ws.Range("E50").Value = "  +0.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.53"
$ws.Range("E51").Value = "  +1.05%  "
